$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end and name it "Final_Matches"
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Final_Matches"

# Header row - reuse the same header style ("s=1": bold, bordered, centered) used by the other sheets
$newSheet.Range("A1").Value = "AZ.CT/LABEL"
$newSheet.Range("B1").Value = "ASCTB.CT/LABEL"
$srcHeader = $wb.Worksheets.Item(1).Range("A1")
$srcHeader.Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @('hematopoeitic stem cell', 'hemopoietic stem cell'),
    @('megakaryocyte-erythroid progenitor cell', 'megakaryocyte-erythroid progenitor cell'),
    @('common lymphoid progenitor', 'common lymphoid progenitor'),
    @('megakaryocyte progenitor cell', 'megakaryocyte progenitor cell'),
    @('monocyte', 'monocyte'),
    @('erythroid lineage cell', 'erythroid lineage cell'),
    @('plasma cell', 'plasma cell'),
    @('memory B cell', 'memory B cell'),
    @('naive B cell', 'naive B cell'),
    @('mature NK cell', 'mature NK T cell'),
    @('transitional stage B cell', 'transitional stage B cell'),
    @('mature natural killer cell', 'mature natural killer cell'),
    @('pro-B cell', 'pro-B cell'),
    @('hematopoeitic multipotent progenitor cell', 'hematopoietic multipotent progenitor cell'),
    @('CD16-negative, CD56-bright natural killer cell, human', 'CD16-negative, CD56-bright natural killer cell'),
    @('lymphocyte of B lineage', 'lymphocyte of B lineage'),
    @('CD14-positive monocyte', 'CD14-positive monocyte'),
    @('dendritic cell', 'dendritic cell, human'),
    @('plasmacytoid dendritic cell', 'plasmacytoid dendritic cell, human'),
    @('CD14-low, CD16-positive monocyte', 'CD14-low, CD16-positive monocyte'),
    @('CD4-positive, alpha-beta T cell', 'T cell'),
    @('CD8-positive, alpha-beta T cell', 'T cell'),
    @('effector CD4-positive, alpha-beta T cell', 'T cell'),
    @('CD4-positive, alpha-beta memory T cell', 'T cell'),
    @('naive thymus-derived CD4-positive, alpha-beta T cell', 'T cell'),
    @('CD8-positive, alpha-beta memory T cell', 'T cell'),
    @('naive thymus-derive CD8-positive, alpha-beta T cell', 'T cell'),
    @('mucosal invariant T cell', 'T cell'),
    @('CD16-positive, CD56-dim natural killer cell, human ', 'mature natural killer cell'),
    @('lymphocyte of B lineage', 'precursor B cell'),
    @('lymphocyte of B lineage', 'immature B cell'),
    @('hematopoietic precursor cell', 'common myeloid progenitor'),
    @('hematopoietic precursor cell', 'promonocyte'),
    @('hematopoietic precursor cell', 'myelocyte'),
    @('hematopoietic precursor cell', 'promyelocyte')
)

$row = 2
foreach ($pair in $data) {
    $newSheet.Cells.Item($row, 1).Value = $pair[0]
    $newSheet.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
